# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K" - strikeouts), replacing the old Strike# counts.
$newK = @{
    2  = 7
    3  = 4
    4  = 7
    5  = 8
    6  = 7
    7  = 9
    8  = 6
    9  = 8
    10 = 6
    11 = 7
    12 = 5
    13 = 5
    14 = 10
    15 = 6
    16 = 10
    17 = 6
    18 = 3
    19 = 10
    20 = 7
    21 = 3
    22 = 4
    23 = 3
    24 = 5
    25 = 8
    26 = 3
    27 = 8
    28 = 5
    29 = 10
    30 = 9
    31 = 8
    32 = 9
    33 = 5
    34 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
